$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Сценарии тестирования")
$ws.Range("C32:C40").Merge()
$ws.Range("C32:C40").HorizontalAlignment = -4108
$ws.Range("C32:C40").VerticalAlignment = -4108
$ws.Range("C32:C40").WrapText = $true
